# Remove 7 client rows (ANDRADE ZAMBRANO MARIUXI ELIZABETH, ARAUZ PARRAGA
# MARIELA PATRICIA, ARBOLEDA ZAMBRANO ROBERTO ANTONIO, ARELLANO CEDEÑO
# DANNY MARCELO, ARMIJOS BARCIA FRACISCO ANTONIO, BAESCORP S.A.S.,
# CHANGKUON AYON JOSE CRISTOBAL) from both the "VENTAS POR GRUPO" and
# "VENTA MENSUAL" sheets, shifting the remaining rows up, then refresh
# the cached summary row that now lands on row 25 of each sheet.

$wb = $excel.ActiveWorkbook
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# Row numbers (in the original layout) of the clients being dropped.
# Row 4 (ARBIZACONSTRUC S.A.) is kept, so it is *not* in this list.
# Deleted highest-to-lowest so earlier deletions don't shift the
# row numbers still queued for removal.
$rowsToDelete = @(9, 8, 7, 6, 5, 3, 2)
foreach ($r in $rowsToDelete) {
    $wsGrupo.Rows.Item($r).Delete()
    $wsMensual.Rows.Item($r).Delete()
}

# The cached "X de 30" counters on the "VENTAS POR GRUPO" summary row
# (now row 25 after the deletions) need their denominator updated to
# reflect the new client count (30 -> 23); the numerators are unaffected
# since none of the removed clients had contributed to them.
$summaryCols = @("C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R")
foreach ($col in $summaryCols) {
    $cell = $wsGrupo.Range($col + "25")
    $cell.Value = ($cell.Value2 -replace "de 30", "de 23")
}

# The cached totals on the "VENTA MENSUAL" summary row (now row 25) are
# plain numbers (not live formulas), so recompute them by hand: they
# drop the contributions of the removed rows (ARELLANO had 86.5 in
# "octubre" and ARAUZ had 489.11 in "noviembre"; the rest were 0).
$wsMensual.Range("C25").Value = 3750.17
$wsMensual.Range("D25").Value = 14407.35
$wsMensual.Range("E25").Value = 96.37
$wsMensual.Range("F25").Value = 206.52
$wsMensual.Range("G25").Value = 0
